$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name) from "Gamma2F-HW30.xpc" to "Gamma2F"
$ws.Name = "Gamma2F"

# Fix a tiny floating point drift in G15 (last-digit rounding change)
$ws.Range("G15").Value = 1.027399527419902

# Copy formatting from the row above so the new row matches existing style (bold/border/centered on col A)
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Append new row 16 of averaged-intensity data (HexGrid-60degTilt5degRes)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.037917103928846
$ws.Range("D16").Value = 0.8587694790858619
$ws.Range("E16").Value = 1.021670918993874
$ws.Range("F16").Value = 1.037917103928846
$ws.Range("G16").Value = 0.9240360028564959
$ws.Range("H16").Value = 1.067212645566715
$ws.Range("I16").Value = 1.02939796776767
$ws.Range("J16").Value = 0.8587694790858619
$ws.Range("K16").Value = 0.9402201990398678
$ws.Range("L16").Value = 0.989068651484357
$ws.Range("M16").Value = 0.9898340196999104
